$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.376.96'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '2.608.94'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E4').Value = '  -0.61%  '
$ws.Range('D5').Value = '''510.97'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = '''154.80'
$ws.Range('E6').Value = '  -1.77%  '
$ws.Range('D7').Value = '''0.997'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '''0.588'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').Value = '2.621.79'
$ws.Range('E9').Value = '  -2.40%  '
$ws.Range('D10').Value = '''6.67'
$ws.Range('E10').Value = '  +2.68%  '
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').Value = '''0.346'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('D14').Value = '3.060.33'
$ws.Range('E14').Value = '  -2.79%  '
$ws.Range('D15').Value = '60.355.98'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '''21.60'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '2.605.83'
$ws.Range('E18').Value = '  -2.95%  '
$ws.Range('D19').Value = '''4.76'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('D20').Value = '''350.55'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = '''10.59'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '''6.14'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '''60.48'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '''0.421'
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('D27').Value = '''0.997'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').Value = '0.0₃0839'
$ws.Range('E28').Value = '  -3.16%  '
$ws.Range('D29').Value = '''7.36'
$ws.Range('E29').Value = '  -2.41%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '''19.43'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').Value = '''151.04'
$ws.Range('E32').Value = '  -4.06%  '
$ws.Range('E33').Value = '  -0.89%  '
$ws.Range('D34').Value = '''5.76'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').Value = '''4.00'
$ws.Range('E35').Value = '  -2.13%  '
$ws.Range('D36').Value = '''1.19'
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('D37').Value = '''0.878'
$ws.Range('E37').Value = '  +4.60%  '
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '''0.843'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '''36.27'
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('D41').Value = '''3.76'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').Value = '''293.99'
$ws.Range('E42').Value = '  -6.58%  '
$ws.Range('D43').Value = '''0.624'
$ws.Range('E43').Value = '  -3.62%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '''0.997'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('D46').Value = '''0.0555'
$ws.Range('E46').Value = '  -3.86%  '
$ws.Range('D47').Value = '''19.73'
$ws.Range('E47').Value = '  -1.63%  '
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').Value = '''0.0234'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = '2.000.53'
$ws.Range('E51').Value = '  -3.61%  '
